$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

$ws.Range("D2").Value = '24.896.20'
$ws.Range("E2").Value = '  +1.75%  '
$ws.Range("D3").Value = '1.669.91'
$ws.Range("E3").Value = '  +0.74%  '
Set-TextValue $ws.Range("D4") '1.003'
$ws.Range("E4").Value = '  +0.01%  '
Set-TextValue $ws.Range("D5") '331.64'
$ws.Range("E5").Value = '  +7.70%  '
Set-TextValue $ws.Range("D6") '1.001'
$ws.Range("E6").Value = '  +0.25%  '
Set-TextValue $ws.Range("D7") '0.3650'
$ws.Range("E7").Value = '  +0.60%  '
Set-TextValue $ws.Range("D8") '46.86'
$ws.Range("E8").Value = '  -1.30%  '
Set-TextValue $ws.Range("D9") '0.3223'
$ws.Range("E9").Value = '  -1.37%  '
Set-TextValue $ws.Range("D10") '1.140'
$ws.Range("E10").Value = '  +1.28%  '
Set-TextValue $ws.Range("D11") '0.07055'
$ws.Range("E11").Value = '  +1.25%  '
$ws.Range("E12").Value = '  +0.15%  '
Set-TextValue $ws.Range("D13") '6.068'
$ws.Range("E13").Value = '  +2.35%  '
Set-TextValue $ws.Range("D14") '19.57'
$ws.Range("D15").Value = '1.664.50'
$ws.Range("E15").Value = '  +0.61%  '
Set-TextValue $ws.Range("D16") '6.619'
$ws.Range("E16").Value = '  +0.00%  '
Set-TextValue $ws.Range("D17") '0.00001046'
$ws.Range("E17").Value = '  +0.24%  '
Set-TextValue $ws.Range("D18") '0.06544'
$ws.Range("E18").Value = '  +0.40%  '
$ws.Range("E19").Value = '  +0.26%  '
Set-TextValue $ws.Range("D20") '78.62'
$ws.Range("E20").Value = '  +2.83%  '
Set-TextValue $ws.Range("D21") '15.82'
$ws.Range("E21").Value = '  +0.68%  '
Set-TextValue $ws.Range("D22") '5.910'
$ws.Range("E22").Value = '  -0.02%  '
Set-TextValue $ws.Range("D23") '12.86'
$ws.Range("E23").Value = '  +2.49%  '
$ws.Range("D24").Value = '24.889.42'
$ws.Range("E24").Value = '  +1.81%  '
Set-TextValue $ws.Range("D25") '2.446'
$ws.Range("E25").Value = '  -0.30%  '
Set-TextValue $ws.Range("D26") '2.395'
$ws.Range("E26").Value = '  +3.74%  '
Set-TextValue $ws.Range("D27") '148.20'
$ws.Range("E27").Value = '  +1.17%  '
Set-TextValue $ws.Range("D28") '18.65'
$ws.Range("E28").Value = '  +1.10%  '
$ws.Range("D29").Value = '1.848.96'
$ws.Range("E29").Value = '  +0.53%  '
Set-TextValue $ws.Range("D30") '125.47'
$ws.Range("E30").Value = '  +0.85%  '
Set-TextValue $ws.Range("D31") '1.172'
$ws.Range("E31").Value = '  -1.74%  '
Set-TextValue $ws.Range("D32") '4.082'
$ws.Range("E32").Value = '  +0.63%  '
Set-TextValue $ws.Range("D33") '5.750'
$ws.Range("E33").Value = '  +2.65%  '
Set-TextValue $ws.Range("D34") '0.08462'
$ws.Range("E34").Value = '  +1.51%  '
Set-TextValue $ws.Range("D35") '1.642'
$ws.Range("E35").Value = '  -2.35%  '
Set-TextValue $ws.Range("D36") '12.25'
$ws.Range("E36").Value = '  -0.88%  '
Set-TextValue $ws.Range("D37") '5.145'
$ws.Range("E37").Value = '  -1.22%  '
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range("D38") '0.06031'
$ws.Range("E38").Value = '  -0.33%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range("D39") '0.02231'
$ws.Range("E39").Value = '  +1.83%  '
$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range("D40") '1.225'
$ws.Range("E40").Value = '  +1.77%  '
Set-TextValue $ws.Range("D41") '0.2084'
$ws.Range("E41").Value = '  +1.66%  '
Set-TextValue $ws.Range("D42") '8.213'
$ws.Range("E42").Value = '  +0.54%  '
Set-TextValue $ws.Range("D43") '1.000'
$ws.Range("E43").Value = '  +0.12%  '
Set-TextValue $ws.Range("D44") '0.5938'
$ws.Range("E44").Value = '  +0.99%  '
Set-TextValue $ws.Range("D45") '13.74'
$ws.Range("E45").Value = '  +8.76%  '
Set-TextValue $ws.Range("D46") '3.863'
$ws.Range("E46").Value = '  +3.48%  '
Set-TextValue $ws.Range("D47") '0.5710'
$ws.Range("E47").Value = '  +2.05%  '
Set-TextValue $ws.Range("D48") '124.81'
$ws.Range("E48").Value = '  +2.35%  '
Set-TextValue $ws.Range("D49") '1.958'
$ws.Range("E49").Value = '  +0.95%  '
Set-TextValue $ws.Range("D50") '0.06991'
$ws.Range("E50").Value = '  +1.19%  '
Set-TextValue $ws.Range("D51") '1.187'
$ws.Range("E51").Value = '  +3.20%  '
